{"js": "// The table's first column holds one benchmark stat per row. This edit:\n//  1. Overwrites the first three summary values with \"0M\".\n//  2. Inserts ten newly-computed rows right after them (200, then nine\n//     tab-separated figures that used to live at the bottom of the table).\n//  3. Collapses the three old tab-separated \"raw dump\" rows at the bottom\n//     down to the single values that used to be at the top of the table.\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// 1) First three rows become \"0M\".\ntable.getCell(0, 0).value = \"0M\";\ntable.getCell(1, 0).value = \"0M\";\ntable.getCell(2, 0).value = \"0M\";\n\n// 2) Insert the ten new rows right after row index 2 (the former \"45\" row).\nconst insertedValues = [\n  [\"200\"],\n  [\"0.00003\"],\n  [\"0.00008\"],\n  [\"0.00005\"],\n  [\"0.00001\"],\n  [\"0.00005\"],\n  [\"0.00005\"],\n  [\"0.00006\"],\n  [\"0.00945\"],\n  [\"100.0\"],\n];\nrows.items[2].insertRows(\"After\", insertedValues.length, insertedValues);\nawait context.sync();\n\n// 3) The three tab-separated rows that used to sit at the end of the table\n// now collapse to single values. Locate them by their distinctive leading\n// number (still present, untouched, after the insert above) rather than by\n// a fixed index so the script is resilient to the row shift.\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst body = table.rows;\nbody.load(\"items\");\nawait context.sync();\n\nconst cells = body.items.map((r) => r.getCell(0));\nfor (const cell of cells) {\n  cell.load(\"value\");\n}\nawait context.sync();\n\nconst replacements = [\n  { prefix: \"97\\t\", value: \"99.98\" },\n  { prefix: \"100\\t\", value: \"0.01\" },\n  { prefix: \"3\\t\", value: \"45\" },\n];\n\nfor (const cell of cells) {\n  const text = cell.value;\n  for (const { prefix, value } of replacements) {\n    if (text && text.indexOf(prefix) === 0) {\n      cell.value = value;\n      break;\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# The table's single column holds one benchmark stat per row. This edit:\n#  1. Overwrites the first three summary values with \"0M\".\n#  2. Inserts ten newly-computed rows right after them (200, then nine\n#     tab-separated figures that used to live at the bottom of the table).\n#  3. Collapses the three old tab-separated \"raw dump\" rows at the bottom\n#     down to the single values that used to be at the top of the table.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# 1) First three rows become \"0M\".\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n\n# 2) Insert the ten new rows right after row 3 (the former \"45\" row), i.e.\n# before the current row 4. Insert in reverse order, always anchored before\n# row 4, so the final on-page order reads 200, 0.00003, ... 100.0.\n$insertedValues = @(\"200\", \"0.00003\", \"0.00008\", \"0.00005\", \"0.00001\", \"0.00005\", \"0.00005\", \"0.00006\", \"0.00945\", \"100.0\")\nfor ($i = $insertedValues.Length - 1; $i -ge 0; $i--) {\n  $newRow = $t.Rows.Add($t.Rows(4))\n  $newRow.Cells(1).Range.Text = $insertedValues[$i]\n}\n\n# 3) The three tab-separated rows that used to sit at the end of the table\n# now collapse to single values. They are the last three rows of the table.\n$lastRow = $t.Rows.Count\n$t.Cell($lastRow - 2, 1).Range.Text = \"99.98\"\n$t.Cell($lastRow - 1, 1).Range.Text = \"0.01\"\n$t.Cell($lastRow, 1).Range.Text = \"45\"\n"}
